$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append extra text to the "Трудовые действия" cell for the first profession row.
$cell = $ws.Range("H2")
$cell.Value = $cell.Text + "; kdshgjkdf"

# Row 2 grew taller to fit the longer wrapped text.
$ws.Rows.Item(2).RowHeight = 57.45

# Selection ended up on G3 after the edit.
$ws.Range("G3").Select()
